$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = [double]"2.0883908691013881E-2"
$ws.Cells.Item(1,2).Value = [double]"1.9334749332489789E-3"
$ws.Cells.Item(1,3).Value = [double]"0.27125785699209837"
$ws.Cells.Item(1,4).Value = [double]"3.4192455792968084E-2"
$ws.Cells.Item(1,5).Value = [double]"5.1947086672103612E-21"
$ws.Cells.Item(1,6).Value = [double]"7.2140060826819202E-2"

$ws.Cells.Item(2,1).Value = [double]"2.0860931471360429E-2"
$ws.Cells.Item(2,2).Value = [double]"7.5525198873936826E-3"
$ws.Cells.Item(2,3).Value = [double]"0.34451965287569397"
$ws.Cells.Item(2,4).Value = [double]"4.044092664500884E-2"
$ws.Cells.Item(2,5).Value = [double]"9.3641097648380672E-21"
$ws.Cells.Item(2,6).Value = [double]"0.70447747311928499"

$ws.Cells.Item(3,1).Value = [double]"2.0815560494486917E-2"
$ws.Cells.Item(3,2).Value = [double]"3.6908230349875248E-3"
$ws.Cells.Item(3,3).Value = [double]"0.31540197987626661"
$ws.Cells.Item(3,4).Value = [double]"3.583766436296492E-2"
$ws.Cells.Item(3,5).Value = [double]"9.5000000014450517E-21"
$ws.Cells.Item(3,6).Value = [double]"0.22105848191344701"

$ws.Cells.Item(4,1).Value = [double]"2.0410650744192693E-2"
$ws.Cells.Item(4,2).Value = [double]"0.15709765247839874"
$ws.Cells.Item(4,3).Value = [double]"0.23469203092355773"
$ws.Cells.Item(4,4).Value = [double]"9.9999999999977801E-2"
$ws.Cells.Item(4,5).Value = [double]"3.5212859309941864E-21"
$ws.Cells.Item(4,6).Value = [double]"9.9999999999999787"

$ws.Cells.Item(5,1).Value = [double]"2.0572461600101691E-2"
$ws.Cells.Item(5,2).Value = [double]"9.1301551138469994E-2"
$ws.Cells.Item(5,3).Value = [double]"0.27427394971343383"
$ws.Cells.Item(5,4).Value = [double]"9.9999999999977801E-2"
$ws.Cells.Item(5,5).Value = [double]"9.5000000014450517E-21"
$ws.Cells.Item(5,6).Value = [double]"9.9999999999999787"

$ws.Cells.Item(6,1).Value = [double]"2.0662036496452894E-2"
$ws.Cells.Item(6,2).Value = [double]"0.23616565190057009"
$ws.Cells.Item(6,3).Value = [double]"0.12037418566801299"
$ws.Cells.Item(6,4).Value = [double]"9.9999999999977801E-2"
$ws.Cells.Item(6,5).Value = [double]"5.8226830347976742E-21"
$ws.Cells.Item(6,6).Value = [double]"9.9999999999999787"

$ws.Cells.Item(7,1).Value = [double]"2.0810355068554765E-2"
$ws.Cells.Item(7,2).Value = [double]"1.0333579257949553E-2"
$ws.Cells.Item(7,3).Value = [double]"0.24851238339370468"
$ws.Cells.Item(7,4).Value = [double]"3.8002889260455808E-2"
$ws.Cells.Item(7,5).Value = [double]"9.5000000014450517E-21"
$ws.Cells.Item(7,6).Value = [double]"8.7686338902047439E-2"

$ws.Cells.Item(8,1).Value = [double]"2.0784198354785472E-2"
$ws.Cells.Item(8,2).Value = [double]"1.8755410504103946E-2"
$ws.Cells.Item(8,3).Value = [double]"0.29524367516283745"
$ws.Cells.Item(8,4).Value = [double]"4.1436946917521485E-2"
$ws.Cells.Item(8,5).Value = [double]"3.7733631137308716E-21"
$ws.Cells.Item(8,6).Value = [double]"0.27852324051423183"

$ws.Cells.Item(9,1).Value = [double]"2.0618377124250426E-2"
$ws.Cells.Item(9,2).Value = [double]"7.5098730231894492E-2"
$ws.Cells.Item(9,3).Value = [double]"0.28347032604155858"
$ws.Cells.Item(9,4).Value = [double]"9.9999999999961162E-2"
$ws.Cells.Item(9,5).Value = [double]"7.378555726882468E-21"
$ws.Cells.Item(9,6).Value = [double]"9.9999999999999609"

$ws.Cells.Item(10,1).Value = [double]"2.0491167196767426E-2"
$ws.Cells.Item(10,2).Value = [double]"7.7194533998224152E-2"
$ws.Cells.Item(10,3).Value = [double]"0.14845809906975785"
$ws.Cells.Item(10,4).Value = [double]"9.9999999999977801E-2"
$ws.Cells.Item(10,5).Value = [double]"7.8124751995675807E-21"
$ws.Cells.Item(10,6).Value = [double]"9.9999999999999787"

$ws.Cells.Item(11,1).Value = [double]"2.0860904450864921E-2"
$ws.Cells.Item(11,2).Value = [double]"1.0151934968936146E-2"
$ws.Cells.Item(11,3).Value = [double]"0.35337322276582422"
$ws.Cells.Item(11,4).Value = [double]"3.7374902727125091E-2"
$ws.Cells.Item(11,5).Value = [double]"5.6356022130339305E-21"
$ws.Cells.Item(11,6).Value = [double]"0.10289723893645383"

$ws.Cells.Item(12,1).Value = [double]"1.97051014418861E-2"
$ws.Cells.Item(12,2).Value = [double]"0.37535935419712763"
$ws.Cells.Item(12,3).Value = [double]"1.2332961612655119E-2"
$ws.Cells.Item(12,4).Value = [double]"9.9999999999977801E-2"
$ws.Cells.Item(12,5).Value = [double]"3.8470918588465785E-21"
$ws.Cells.Item(12,6).Value = [double]"9.9999999999999787"

$ws.Cells.Item(13,1).Value = [double]"2.0848532399393927E-2"
$ws.Cells.Item(13,2).Value = [double]"1.048497754959485E-2"
$ws.Cells.Item(13,3).Value = [double]"0.34120252311446819"
$ws.Cells.Item(13,4).Value = [double]"3.9423978545248983E-2"
$ws.Cells.Item(13,5).Value = [double]"9.5000000014450517E-21"
$ws.Cells.Item(13,6).Value = [double]"0.26759564881345532"

$ws.Cells.Item(14,1).Value = [double]"2.0914237207588504E-2"
$ws.Cells.Item(14,2).Value = [double]"4.8797269025110156E-3"
$ws.Cells.Item(14,3).Value = [double]"0.28993868501251596"
$ws.Cells.Item(14,4).Value = [double]"3.454064867966504E-2"
$ws.Cells.Item(14,5).Value = [double]"9.5136863374092964E-21"
$ws.Cells.Item(14,6).Value = [double]"5.0139460186337859E-2"

$ws.Cells.Item(15,1).Value = [double]"2.1053144130016518E-2"
$ws.Cells.Item(15,2).Value = [double]"5.1902395790734322E-2"
$ws.Cells.Item(15,3).Value = [double]"0.37910896331036353"
$ws.Cells.Item(15,4).Value = [double]"5.256957035208408E-2"
$ws.Cells.Item(15,5).Value = [double]"9.5038933455566795E-21"
$ws.Cells.Item(15,6).Value = [double]"1.5658298358075942"

$ws.Cells.Item(16,1).Value = [double]"2.0918251657774008E-2"
$ws.Cells.Item(16,2).Value = [double]"2.1291626147412455E-2"
$ws.Cells.Item(16,3).Value = [double]"0.34457000483966999"
$ws.Cells.Item(16,4).Value = [double]"5.7209863257016394E-2"
$ws.Cells.Item(16,5).Value = [double]"3.6110569609447987E-21"
$ws.Cells.Item(16,6).Value = [double]"2.7602987209048933"

